# Insert two new data rows right before the current row 576, pushing the
# existing rows 576..687 down to 578..689 (dimension grows from A1:T687 to
# A1:T689), then populate the two freshly-inserted rows with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("576:577").Insert()

# New row 576: Femacal de La Calera / Coquimbo / Naranja / Valencia / Primera
$ws.Range("A576").Value = 3
$ws.Range("B576").Value = "Femacal de La Calera"
$ws.Range("C576").Value = "Coquimbo"
$ws.Range("D576").Value = 44637
$ws.Range("E576").Value = 5
$ws.Range("F576").Value = "Fruta"
$ws.Range("G576").Value = 100102
$ws.Range("H576").Value = "Cítricos"
$ws.Range("I576").Value = 100102005
$ws.Range("J576").Value = "Naranja"
$ws.Range("K576").Value = "Valencia"
$ws.Range("L576").Value = "Primera"
$ws.Range("M576").Value = 138
$ws.Range("N576").Value = 7500
$ws.Range("O576").Value = 8000
$ws.Range("P576").Value = 7754
$ws.Range("Q576").Value = "`$/malla 13 kilos"
$ws.Range("R576").Value = "Provincia de Quillota"
$ws.Range("S576").Value = 596
$ws.Range("T576").Value = 13

# New row 577: Femacal de La Calera / Coquimbo / Naranja / Valencia / Segunda
$ws.Range("A577").Value = 3
$ws.Range("B577").Value = "Femacal de La Calera"
$ws.Range("C577").Value = "Coquimbo"
$ws.Range("D577").Value = 44637
$ws.Range("E577").Value = 5
$ws.Range("F577").Value = "Fruta"
$ws.Range("G577").Value = 100102
$ws.Range("H577").Value = "Cítricos"
$ws.Range("I577").Value = 100102005
$ws.Range("J577").Value = "Naranja"
$ws.Range("K577").Value = "Valencia"
$ws.Range("L577").Value = "Segunda"
$ws.Range("M577").Value = 80
$ws.Range("N577").Value = 6000
$ws.Range("O577").Value = 6000
$ws.Range("P577").Value = 6000
$ws.Range("Q577").Value = "`$/malla 13 kilos"
$ws.Range("R577").Value = "Provincia de Quillota"
$ws.Range("S577").Value = 462
$ws.Range("T577").Value = 13
